$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flatten the multi-line instructional comment strings (remove embedded
# line breaks, turning them into single-line guidance) for the relevant
# template cells, and add the new import instructions text.
$ws.Range("E2").Value = "Enter either: Female Male"
$ws.Range("F2").Value = "YYYY-MM-DD (formatted as text) ex. 2020-08-01"
$ws.Range("H2").Value = "Format per country 919-555-1212"
$ws.Range("I2").Value = "Enter either: Home   Work  Mobile"
$ws.Range("K2").Value = "Enter either: Home   Work  Mobile"

# Reset column K (11) back to (near) the default width.
$ws.Columns(11).ColumnWidth = 10.166666666666666

# Move the active selection to F5, as left by the author after editing.
$ws.Range("F5").Select()
